$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sorted / cleaned weather-condition list, replacing the old header +
# trailing-space-padded values. "condicao_metereologica" header is gone and
# "n/a" has been inserted in alphabetical order.
$values = @("Céu Claro", "Chuva", "Garoa/Chuvisco", "Granizo", "Ignorado", "n/a", "Neve", "Nevoeiro/Neblina", "Nublado", "Sol", "Vento")

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# A1 was the bold header style; it's now a regular value like the rest.
$ws.Range("A1").Font.Bold = $false

# Reflect the saved selection/active cell.
$ws.Range("A3").Select()
